$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest quotation row (2025-09-16) below the existing data,
# mirroring the date cell format of column A and the comma-decimal text
# values used in columns B:E for every prior row.
$ws.Range("A12").Value = 45916
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat

$ws.Range("B12").Value = "21,2666"
$ws.Range("C12").Value = "15,0221"
$ws.Range("D12").Value = "15,0023"
$ws.Range("E12").Value = "15,0023"
